$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 15627200
$ws.Range("J113").Value = 31251924
$ws.Range("L113").Value = 31251924
$ws.Range("N113").Value = -31258432

$ws.Range("H125").Value = 2602.4546
$ws.Range("I125").Value = 2762.7
$ws.Range("K125").Value = 24864.3
$ws.Range("M125").Value = -22404.3

$ws.Range("H127").Value = 942.95
$ws.Range("I127").Value = 604.53845
$ws.Range("J127").Value = 1571.4286
$ws.Range("K127").Value = 1813.61535
$ws.Range("L127").Value = 4714.2858
$ws.Range("M127").Value = 3146.38465
$ws.Range("N127").Value = -14634.2858

$ws.Range("H129").Value = 1249.4348
$ws.Range("I129").Value = 851.4286
$ws.Range("J129").Value = 1423.5625
$ws.Range("K129").Value = 2554.2858
$ws.Range("L129").Value = 4270.6875
$ws.Range("M129").Value = 2445.7142
$ws.Range("N129").Value = -14270.6875

$ws.Range("H136").Value = 49325
$ws.Range("J136").Value = 49325
$ws.Range("L136").Value = 49325
$ws.Range("N136").Value = -59525

$ws.Range("H139").Value = 50270
$ws.Range("J139").Value = 50270
$ws.Range("L139").Value = 50270
$ws.Range("N139").Value = -60550

$ws.Range("H140").Value = 74666.664
$ws.Range("J140").Value = 74666.664
$ws.Range("L140").Value = 74666.664
$ws.Range("N140").Value = -85026.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11423.25
$ws.Range("I32").Value = 9597.392
$ws.Range("J32").Value = 16619.924
$ws.Range("K32").Value = 9597.392
$ws.Range("L32").Value = 16619.924
$ws.Range("M32").Value = -9310.392
$ws.Range("N32").Value = -17193.924

$ws.Range("H61").Value = 1114888.6
$ws.Range("I61").Value = 3697.25
$ws.Range("J61").Value = 2003841.8
$ws.Range("K61").Value = 3697.25
$ws.Range("L61").Value = 2003841.8
$ws.Range("M61").Value = -3485.25
$ws.Range("N61").Value = -2004265.8

$ws.Range("H122").Value = 3213992
$ws.Range("I122").Value = 4283656
$ws.Range("K122").Value = 12850968
$ws.Range("M122").Value = -12848518

$ws.Range("H136").Value = 1114888.6
$ws.Range("I136").Value = 3697.25
$ws.Range("J136").Value = 2003841.8
$ws.Range("K136").Value = 11091.75
$ws.Range("L136").Value = 6011525.4
$ws.Range("M136").Value = -8541.75
$ws.Range("N136").Value = -6016625.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 429603.16
$ws.Range("I107").Value = 601064.4
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 601064.4
$ws.Range("L107").Value = 950
$ws.Range("M107").Value = -599144.4
$ws.Range("N107").Value = -4790

$ws.Range("H137").Value = 59605
$ws.Range("J137").Value = 59605
$ws.Range("L137").Value = 59605
$ws.Range("N137").Value = -69805

$ws.Range("H138").Value = 60780
$ws.Range("J138").Value = 60780
$ws.Range("L138").Value = 60780
$ws.Range("N138").Value = -71060

$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10646393
$ws.Range("I31").Value = 2296.4546
$ws.Range("J31").Value = 13898756
$ws.Range("K31").Value = 2296.4546
$ws.Range("L31").Value = 13898756
$ws.Range("M31").Value = -2001.4546
$ws.Range("N31").Value = -13899346

$ws.Range("H34").Value = 10646393
$ws.Range("I34").Value = 2296.4546
$ws.Range("J34").Value = 13898756
$ws.Range("K34").Value = 2296.4546
$ws.Range("L34").Value = 13898756
$ws.Range("M34").Value = -2094.4546
$ws.Range("N34").Value = -13899160

$ws.Range("H135").Value = 37547.31
$ws.Range("J135").Value = 37547.31
$ws.Range("L135").Value = 37547.31
$ws.Range("N135").Value = -47687.31

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 5000455
$ws.Range("I98").Value = 305
$ws.Range("J98").Value = 6250492.5
$ws.Range("K98").Value = 915
$ws.Range("L98").Value = 18751477.5
$ws.Range("M98").Value = 583
$ws.Range("N98").Value = -18754473.5

$ws.Range("H107").Value = 395.93103
$ws.Range("I107").Value = 335.0909
$ws.Range("J107").Value = 587.1429000000001
$ws.Range("K107").Value = 1005.2727
$ws.Range("L107").Value = 1761.4287
$ws.Range("M107").Value = 914.7273
$ws.Range("N107").Value = -5601.4287

$ws.Range("H121").Value = 777.25
$ws.Range("J121").Value = 1031.1875
$ws.Range("L121").Value = 3093.5625
$ws.Range("N121").Value = -5713.5625

$ws.Range("H131").Value = 2128846
$ws.Range("J131").Value = 1457.9678
$ws.Range("L131").Value = 4373.903399999999
$ws.Range("N131").Value = -14453.9034

$ws.Range("H139").Value = 4689.7144
$ws.Range("I139").Value = 5816.1904
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 17448.5712
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -12308.5712
$ws.Range("N139").Value = -19280

$ws.Range("H140").Value = 8370.111000000001
$ws.Range("I140").Value = 8370.111000000001
$ws.Range("K140").Value = 25110.333
$ws.Range("M140").Value = -19930.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3212.4666
$ws.Range("I102").Value = 2715.8572
$ws.Range("K102").Value = 2715.8572
$ws.Range("M102").Value = -1093.8572

$ws.Range("H131").Value = 33333
$ws.Range("J131").Value = 33333
$ws.Range("L131").Value = 33333
$ws.Range("N131").Value = -43413

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3499.889
$ws.Range("I61").Value = 3401
$ws.Range("J61").Value = 3579
$ws.Range("K61").Value = 3401
$ws.Range("L61").Value = 3579
$ws.Range("M61").Value = -3199
$ws.Range("N61").Value = -3983

$ws.Range("H113").Value = 3499.889
$ws.Range("I113").Value = 3401
$ws.Range("J113").Value = 3579
$ws.Range("K113").Value = 3401
$ws.Range("L113").Value = 3579
$ws.Range("M113").Value = -1231
$ws.Range("N113").Value = -7919

$ws.Range("H141").Value = 82958.42999999999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 82958.42999999999
$ws.Range("K141").Value = 0
$ws.Range("N141").Value = -93318.42999999999
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6176792
$ws.Range("I136").Value = 5752.4546
$ws.Range("J136").Value = 10419381
$ws.Range("K136").Value = 17257.3638
$ws.Range("L136").Value = 31258143
$ws.Range("M136").Value = -14707.3638
$ws.Range("N136").Value = -31263243
